$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Configuration")
$wsAppEnv = $wb.Worksheets.Item("Application_Environment")

# Configuration!B3: Y -> N
$wsConfig.Range("B3").Value = "N"

# Clear row 16 on Application_Environment (RGT_AUSCountry row) and remove its yellow fill
$rngRow16 = $wsAppEnv.Range("A16:H16")
$rngRow16.ClearContents()
$rngRow16.Interior.Pattern = -4142

$wsAppEnv.Range("I16").Clear()

# Selection / active-cell bookkeeping
$wsConfig.Range("B8").Select()
$wsAppEnv.Range("F19").Select()
$wsAppEnv.Activate()
